$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the assignment to be stored as text (columns D/E hold
    # formatted strings like "26.607.20" or "  +0.98%  ", not numbers),
    # then strip the temporary "@" number format so the cell keeps the
    # original (default) style, matching the source file.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextCell $ws.Range("D2") "26.595.43"
Set-TextCell $ws.Range("E2") "  +0.96%  "
Set-TextCell $ws.Range("D3") "1.819.69"
Set-TextCell $ws.Range("E3") "  +1.69%  "
Set-TextCell $ws.Range("D4") "1.007"
Set-TextCell $ws.Range("E4") "  -0.10%  "
Set-TextCell $ws.Range("D5") "1.006"
Set-TextCell $ws.Range("E5") "  -0.07%  "
Set-TextCell $ws.Range("D6") "305.61"
Set-TextCell $ws.Range("E6") "  -0.36%  "
Set-TextCell $ws.Range("D7") "0.4666"
Set-TextCell $ws.Range("E7") "  +2.29%  "
Set-TextCell $ws.Range("D8") "0.3598"
Set-TextCell $ws.Range("E8") "  -0.95%  "
Set-TextCell $ws.Range("D9") "46.13"
Set-TextCell $ws.Range("E9") "  -0.19%  "
Set-TextCell $ws.Range("D10") "0.07123"
Set-TextCell $ws.Range("E10") "  +0.62%  "
Set-TextCell $ws.Range("D11") "0.9018"
Set-TextCell $ws.Range("E11") "  +3.15%  "
Set-TextCell $ws.Range("D12") "0.07780"
Set-TextCell $ws.Range("E12") "  -0.33%  "
Set-TextCell $ws.Range("D13") "19.43"
Set-TextCell $ws.Range("E13") "  -0.23%  "
Set-TextCell $ws.Range("D14") "1.877.60"
Set-TextCell $ws.Range("E14") "  +4.17%  "
Set-TextCell $ws.Range("D15") "5.245"
Set-TextCell $ws.Range("E15") "  -0.43%  "
Set-TextCell $ws.Range("D16") "6.326"
Set-TextCell $ws.Range("E16") "  +0.23%  "
Set-TextCell $ws.Range("D17") "87.44"
Set-TextCell $ws.Range("E18") "  -0.13%  "
Set-TextCell $ws.Range("D19") "0.000008556"
Set-TextCell $ws.Range("E19") "  +0.47%  "
Set-TextCell $ws.Range("E20") "  -0.05%  "
Set-TextCell $ws.Range("D21") "26.643.01"
Set-TextCell $ws.Range("E21") "  +1.02%  "
Set-TextCell $ws.Range("D22") "14.20"
Set-TextCell $ws.Range("E22") "  -0.25%  "
Set-TextCell $ws.Range("D23") "5.009"
Set-TextCell $ws.Range("E23") "  +0.60%  "
Set-TextCell $ws.Range("D24") "10.56"
Set-TextCell $ws.Range("E24") "  +0.63%  "
Set-TextCell $ws.Range("D25") "1.932"
Set-TextCell $ws.Range("E25") "  -2.36%  "
Set-TextCell $ws.Range("D26") "151.97"
Set-TextCell $ws.Range("E26") "  -0.06%  "
Set-TextCell $ws.Range("D27") "17.88"
Set-TextCell $ws.Range("E27") "  -0.01%  "
Set-TextCell $ws.Range("D28") "1.982"
Set-TextCell $ws.Range("E28") "  -2.45%  "
Set-TextCell $ws.Range("D29") "113.52"
Set-TextCell $ws.Range("E29") "  +1.19%  "
Set-TextCell $ws.Range("D30") "4.815"
Set-TextCell $ws.Range("E30") "  -0.67%  "
Set-TextCell $ws.Range("D31") "0.08775"
Set-TextCell $ws.Range("E31") "  +1.27%  "
Set-TextCell $ws.Range("D32") "3.136"
Set-TextCell $ws.Range("E32") "  +2.86%  "
Set-TextCell $ws.Range("D33") "2.783"
Set-TextCell $ws.Range("E33") "  +4.81%  "
Set-TextCell $ws.Range("D34") "0.7334"
Set-TextCell $ws.Range("E34") "  +1.98%  "
Set-TextCell $ws.Range("D35") "4.437"
Set-TextCell $ws.Range("E35") "  -0.06%  "
Set-TextCell $ws.Range("E36") "  +1.79%  "
Set-TextCell $ws.Range("D37") "1.076"
Set-TextCell $ws.Range("E37") "  -0.08%  "
Set-TextCell $ws.Range("D38") "0.01932"
Set-TextCell $ws.Range("E38") "  -0.58%  "
Set-TextCell $ws.Range("D39") "2.911"
Set-TextCell $ws.Range("E39") "  +1.48%  "
Set-TextCell $ws.Range("D40") "0.05118"
Set-TextCell $ws.Range("E40") "  +0.41%  "
Set-TextCell $ws.Range("D41") "0.5064"
Set-TextCell $ws.Range("E41") "  -0.52%  "
Set-TextCell $ws.Range("E42") "  -1.47%  "
Set-TextCell $ws.Range("D43") "0.1496"
Set-TextCell $ws.Range("E43") "  -1.36%  "
Set-TextCell $ws.Range("D44") "7.995"
Set-TextCell $ws.Range("E44") "  -0.20%  "
Set-TextCell $ws.Range("D45") "0.4681"
Set-TextCell $ws.Range("E45") "  +0.65%  "
Set-TextCell $ws.Range("D46") "1.006"
Set-TextCell $ws.Range("E46") "  -0.11%  "
Set-TextCell $ws.Range("D47") "10.00"
Set-TextCell $ws.Range("E47") "  +1.42%  "
Set-TextCell $ws.Range("D48") "98.52"
Set-TextCell $ws.Range("E48") "  -1.30%  "
Set-TextCell $ws.Range("D49") "1.564"
Set-TextCell $ws.Range("E49") "  -1.35%  "
Set-TextCell $ws.Range("D50") "0.06004"
Set-TextCell $ws.Range("E50") "  +0.56%  "
Set-TextCell $ws.Range("D51") "63.66"
Set-TextCell $ws.Range("E51") "  -0.33%  "
